$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy formatting (bold font,
# border, center/top alignment) from the existing "IP" header cell (H1)
# so the new headers look consistent with the rest of the header row.
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 10).PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Fill in the new I/J columns for every data row: I is always 1,
# J mirrors the existing IP value in column H.
for ($r = 2; $r -le 37; $r++) {
    $h = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $h
}
